$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rebuild rows 2-27 (header row 1 stays the same) with the revised facility list / heuristic
$data = @(
    @(0, 0, "INNOCENTI SRL", 41.955471, 12.7640919, "facility"),
    @(1, 10, "FITALS S.R.L.", 41.9404295, 12.632209, "facility"),
    @(2, 11, "PORCARELLI GINO & C. SRL", 41.8863388, 12.6937481, "facility"),
    @(3, 57, "MATERIALE COMMERCIO RICUPERO MCR SRL", 41.8793057, 12.5530678, "facility"),
    @(4, 68, "GE.SER. SRL", 41.9507218, 12.7554548, "facility"),
    @(5, 81, "RECYCLING PFU SRL", 41.2759554, 13.7880158, "facility"),
    @(6, 139, "RICREA SRL", 41.7082726, 12.5870176, "facility"),
    @(7, 205, "SIECO S.R.L.", 42.4539311, 12.0871977, "facility"),
    @(8, 211, "RIME 1", 41.8345359, 12.4310169, "facility"),
    @(9, 270, "ECOLOGIC SYSTEM COMPANY SRL", 41.6942634, 12.534113, "facility"),
    @(10, 349, "CE.STRA SRL", 41.8479537, 12.6346416, "facility"),
    @(11, 441, "FEROCART SRL", 41.89457, 12.37798, "facility"),
    @(12, 484, "SIDERURGICA TIBURTINA S.R.L.", 41.9602459, 12.7015494, "facility"),
    @(13, 524, "SO.GE.RI.T SRL", 41.4132753, 13.1695125, "facility"),
    @(15, 663, "N.I.ECO S.P.A.", 41.9271213, 12.6238636, "facility"),
    @(16, 772, "RICICLA CENTRO ITALIA SRL", 41.8842861, 12.704119, "facility"),
    @(18, 965, "MEDIASERVICE RECYCLING SRL", 41.9319972, 12.6140927, "facility"),
    @(22, 1368, "ECO LOGICA 2000 SRL", 41.7824934, 12.604908, "facility"),
    @(25, 2134, "ECO LOGICA 2000 SRL", 41.7931483, 12.5359218, "facility"),
    @(27, 3258, "IBIOS SRL", 42.0461171, 12.5552946, "facility"),
    @(28, 3403, "EASY TONER SNC", 41.9743277, 12.5075322, "facility"),
    @(30, 4134, "SE.IN SRL", 41.64394920000001, 13.2197606, "facility"),
    @(32, 4763, "CARLUCCIO ROTTAMI DI ROCCA & C.SNC", 41.955343, 12.765942, "facility"),
    @(33, 4881, "FITALS SRL", 41.9404295, 12.632209, "facility"),
    @(34, 5463, "ECOSYSTEM SPA", 41.7064274, 12.5457223, "facility"),
    @(35, 6826, "TECNOGARDEN SERVICE SRL", 41.8916652, 12.6232211, "facility"),
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $row++
}

# Newly added rows (12-27) need the same formatting already applied to column A
# (bold, thin border, centered/top-aligned) - copy it from an existing styled cell
# instead of rebuilding it, so the style is reused rather than duplicated.
$ws.Range("A2").Copy()
$ws.Range("A12:A27").PasteSpecial(-4122)
$excel.CutCopyMode = $false
